$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 6
$ws.Range("G2").Value = -3
$ws.Range("H2").Value = 13

$ws.Range("D2").Select()
